$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '61.508.81'
$ws.Range('E2').Value = '  +2.08%  '
$ws.Range('D3').Value = '2.658.68'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.37'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.03'
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.62'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('E10').Value = '  +4.51%  '
$ws.Range('E11').Value = '  +3.07%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '3.125.71'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('E14').Value = '  +7.21%  '
$ws.Range('D15').Value = '61.388.45'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('E16').Value = '  +4.25%  '
$ws.Range('D17').Value = '2.661.85'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.69'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('E19').Value = '  +3.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '356.49'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.92'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.73'
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('E25').Value = '  +3.39%  '
$ws.Range('E26').Value = '  +5.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.01'
$ws.Range('E28').Value = '  +8.26%  '
$ws.Range('D29').Value = '0.0₃0828'
$ws.Range('E29').Value = '  +4.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.97'
$ws.Range('E30').Value = '  +8.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.40'
$ws.Range('E31').Value = '  +2.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.997'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.19'
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.14'
$ws.Range('E34').Value = '  +15.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.69'
$ws.Range('E35').Value = '  +8.58%  '
$ws.Range('E36').Value = '  +7.51%  '
$ws.Range('E37').Value = '  +19.35%  '
$ws.Range('E38').Value = '  +5.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '347.23'
$ws.Range('E39').Value = '  +10.52%  '
$ws.Range('E40').Value = '  +6.23%  '
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.42'
$ws.Range('E42').Value = '  +7.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.89'
$ws.Range('E43').Value = '  +5.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.30'
$ws.Range('E44').Value = '  +6.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0580'
$ws.Range('E45').Value = '  +5.22%  '
$ws.Range('E46').Value = '  +4.00%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.81'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0254'
$ws.Range('E48').Value = '  +5.24%  '
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.996'
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').Value = '2.103.43'
$ws.Range('E51').Value = '  +3.51%  '
